# Generate Report for Archive
# 1. Update status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview sheet zh-cn/de-de status columns, and the Status
#    column on each per-locale sheet).
# 2. Shrink the "Status" related columns (Overview!E:F and the Status column
#    on each locale sheet) to their new, narrower width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- 1. Replace "Ready for handoff" with "In Translation" -------------------
# Overview sheet: columns E (zh-cn) and F (de-de), data rows 2-4
$overview.Range("E2:F4").Value = "In Translation"

# Per-locale sheets: column C ("Status"), data rows 2-4
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# --- 2. Narrow the status columns --------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
